$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Invalid username or password"
$ws.Range("D3").Value = "Invalid username or password"

$ws.Range("F3").Select()
